$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

# Each entry: row (1-indexed within the table), col (1-indexed), new text
# (Row/col verified against the original document's cell contents.)
$cellEdits = @(
    @(1, 1, "38÷9="),
    @(1, 2, "23÷8="),
    @(1, 3, "44÷4="),
    @(1, 4, "65÷2="),
    @(1, 5, "18÷5="),
    @(5, 1, "88÷7="),
    @(5, 2, "31÷4="),
    @(5, 3, "46÷4="),
    @(5, 4, "30÷5="),
    @(5, 5, "29÷9="),
    @(9, 1, "73÷9="),
    @(9, 2, "77÷3="),
    @(9, 3, "73÷7="),
    @(9, 4, "87÷8="),
    @(9, 5, "82÷7="),
    @(13, 1, "78÷6="),
    @(13, 2, "20÷8="),
    @(13, 3, "74÷6="),
    @(13, 4, "90÷2="),
    @(13, 5, "76÷2="),
    @(17, 1, "63÷4="),
    @(17, 2, "29÷4="),
    @(17, 3, "97÷4="),
    @(17, 4, "40÷6="),
    @(17, 5, "81÷5=")
)

foreach ($edit in $cellEdits) {
    $row = $edit[0]
    $col = $edit[1]
    $new = $edit[2]

    $cell = $table.Cell($row, $col)
    $range = $cell.Range
    $range.Text = $new
}
